$wb = $excel.ActiveWorkbook

# --- classFields sheet: reorder field rows within each class block ---
$ws = $wb.Worksheets.Item("classFields")

$ws.Range("B2").Value = "repository"
$ws.Range("D2").Value = "com.zatribune.spring.ecommerce.payments.db.repository.CustomerRepository"

$ws.Range("B3").Value = "log"
$ws.Range("D3").Value = "org.slf4j.Logger"

$ws.Range("B4").Value = "log"
$ws.Range("D4").Value = "org.slf4j.Logger"

$ws.Range("B5").Value = "orderService"
$ws.Range("D5").Value = "com.zatribune.spring.ecommerce.payments.service.OrderService"

$ws.Range("B8").Value = "id"
$ws.Range("D8").Value = "java.lang.Long"

$ws.Range("B9").Value = "amountReserved"
$ws.Range("D9").Value = "int"

$ws.Range("B10").Value = "amountReserved"
$ws.Range("D10").Value = "int"

$ws.Range("B11").Value = "id"
$ws.Range("D11").Value = "java.lang.Long"

$ws.Range("B12").Value = "amountAvailable"
$ws.Range("D12").Value = "int"

$ws.Range("B13").Value = "name"
$ws.Range("D13").Value = "java.lang.String"

$ws.Range("B14").Value = "template"
$ws.Range("D14").Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Range("B15").Value = "SOURCE"
$ws.Range("D15").Value = "domain.OrderSource"

$ws.Range("B16").Value = "log"
$ws.Range("D16").Value = "org.slf4j.Logger"

# --- classNumberOfLines sheet: CustomerRepository line count fix ---
$ws2 = $wb.Worksheets.Item("classNumberOfLines")
$ws2.Range("B4").Value = "0"
